# Append a new registration row (row 8) to the "Registrations" sheet,
# mirroring the existing rows above it.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 8

$ws.Range("A$row").Value2 = "AGNIVA BHATTACHARJEE"
$ws.Range("B$row").Value2 = "bhattacharjee.agniva.jobs@gmail.com"

# Phone + Pass-Out-Year must stay text (so a leading "0" on the phone
# number, or a bare "2020" year, aren't coerced into numbers). Flip the
# cell to a text number-format just long enough to type the value in as
# a string, then flip it back to the workbook's default ("Normal") style
# so the cell does not end up tagged with a different style index than
# its neighbours.
$ws.Range("C$row").NumberFormat = "@"
$ws.Range("C$row").Value2 = "08420880979"
$ws.Range("C$row").Style = "Normal"

$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value2 = "2020"
$ws.Range("D$row").Style = "Normal"

$ws.Range("E$row").Value2 = "3ba20b57-9c39-4880-8ae1-203fa0db0544"

# "Entered" column is blank for every row (stored as an empty string in
# the source file); leave it unset so the new cell is likewise empty.
$ws.Range("F$row").Value2 = ""
